$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update question text for rows 3 and 4 (column C) to add numbering prefixes
$ws.Range("C3").Value = "1.1 Um welche Unterlage handelt es sich"
$ws.Range("C4").Value = "1.1.1 Ist diese Unterlage vollständig?"

# Apply General number format to A2:B13 (creates a new cell style entry)
$ws.Range("A2:B13").NumberFormat = "General"

# Update the selection to A2:B13
$ws.Range("A2:B13").Select
